$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the new J1 / J2 designators, right after the D1 row (row 6),
# pushing L1 and everything below it down by two rows.
$ws.Rows("7:8").Insert()

# Populate the two newly inserted rows.
$ws.Range("A7").Value = "J1"
$ws.Range("B7").Value = " 52.3000mm"
$ws.Range("C7").Value = "  -60.0000mm"
$ws.Range("D7").Value = " top"
$ws.Range("E7").Value = 90

$ws.Range("A8").Value = "J2"
$ws.Range("B8").Value = " 73.1000mm"
$ws.Range("C8").Value = "  -60.0000mm"
$ws.Range("D8").Value = " top"
$ws.Range("E8").Value = 270

# The "Mid Y" column (C) gained an extra leading space in front of the minus sign
# for every existing data row. Rewrite each one.
$ws.Range("C2").Value  = "  -54.6000mm"   # C1
$ws.Range("C3").Value  = "  -66.8300mm"   # C2
$ws.Range("C4").Value  = "  -67.4500mm"   # C3
$ws.Range("C5").Value  = "  -56.5000mm"   # C4
$ws.Range("C6").Value  = "  -60.0000mm"   # D1
$ws.Range("C9").Value  = "  -52.4500mm"   # L1
$ws.Range("C10").Value = "  -62.1000mm"   # Q1
$ws.Range("C11").Value = "  -52.9500mm"   # Q2
$ws.Range("C12").Value = "  -65.7000mm"   # Q3
$ws.Range("C13").Value = "  -64.9000mm"   # R1
$ws.Range("C14").Value = "  -64.8000mm"   # R2
$ws.Range("C15").Value = "  -61.7000mm"   # R3
$ws.Range("C16").Value = "  -63.0000mm"   # R4
$ws.Range("C17").Value = "  -68.4000mm"   # R5
$ws.Range("C18").Value = "  -56.0000mm"   # R6
$ws.Range("C19").Value = "  -67.0500mm"   # U1
$ws.Range("C20").Value = "  -52.9500mm"   # U2
$ws.Range("C21").Value = "  -67.3500mm"   # U3

# A handful of rotation values were also tweaked for manufacturing.
$ws.Range("E10").Value = 0     # Q1: 180 -> 0
$ws.Range("E11").Value = 90    # Q2: 270 -> 90
$ws.Range("E12").Value = 90    # Q3: 270 -> 90
$ws.Range("E19").Value = 180   # U1: 0   -> 180
$ws.Range("E20").Value = 270   # U2: 90  -> 270
$ws.Range("E21").Value = 0     # U3: 180 -> 0

# The active selection moved too.
$ws.Range("F22").Select()
